# SCD0338-001 fix: update "Test Daily Activity" note + schedule times for
# rows 2 and 3, and move the active selection from Q3 to P3.
#
# Row 2 (RUN=43988): "Test Daily Activity 2" -> "Test Daily Activity 3"
#                      03:50 PM -> 07:05 PM ; 04:30 PM -> 08:30 PM
# Row 3 (RUN=103258): "Test Daily Activity 2" -> "Test Daily Activity 3"
#                      03:50 PM -> 08:50 PM ; 04:30 PM -> 10:00 PM
#
# The O/P time strings get a leading apostrophe so they stay literal text
# (matching the sheet's existing quotePrefix text style for those columns)
# instead of Excel reinterpreting "HH:MM PM" as a time value/serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L2").Value = "Test Daily Activity 3"
$ws.Range("O2").Value = "'07:05 PM"
$ws.Range("P2").Value = "'08:30 PM"

$ws.Range("L3").Value = "Test Daily Activity 3"
$ws.Range("O3").Value = "'08:50 PM"
$ws.Range("P3").Value = "'10:00 PM"

# The cached TODAY() display (N2/N3) recalculates automatically from the
# unchanged formula; no explicit write needed.

# Move the selection from Q3 to P3, matching the saved sheet view.
$ws.Range("P3").Select()
